{"js": "// Adjust the workshop start time from \"ab 19:30\" to \"ab 19:15\" and move the\n// \"_GoBack\" bookmark so it still marks the position of the last edit (right\n// after the time text, before the closing parenthesis) - mirroring what Word\n// itself does when you retype text at the location the bookmark used to sit.\n\nconst body = context.document.body;\n\n// 1) Remove the (now stale) \"_GoBack\" bookmark that currently sits inside the\n//    \"Exploratives/Manuelles Testen\" bullet.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the run containing the old start time and replace its text only -\n//    this keeps the surrounding \"(\" and \")\" runs untouched, exactly like the\n//    diff shows.\nconst timeResults = body.search(\"ab 19:30\", { matchCase: true, matchWholeWord: false });\ntimeResults.load(\"text\");\nawait context.sync();\n\nif (timeResults.items.length === 0) {\n  throw new Error('Could not find \"ab 19:30\" in the document.');\n}\n\nconst timeRange = timeResults.items[0];\ntimeRange.insertText(\"ab 19:15\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-add the \"_GoBack\" bookmark right after the new time text (i.e. right\n//    before the closing \")\").\nconst newTimeResults = body.search(\"ab 19:15\", { matchCase: true, matchWholeWord: false });\nnewTimeResults.load(\"text\");\nawait context.sync();\n\nconst newTimeRange = newTimeResults.items[0];\nconst afterTime = newTimeRange.getRange(Word.RangeLocation.end);\nafterTime.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Adjust the workshop start time from \"ab 19:30\" to \"ab 19:15\" and move the\n# \"_GoBack\" bookmark so that it keeps marking the spot of the most recent\n# edit - i.e. right after the new time text, just before the closing \")\" -\n# instead of its old, now-stale position inside \"Exploratives/Manuelles\n# Testen\".\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old \"_GoBack\" bookmark (it currently sits inside the\n#    \"Exploratives/Manuelles Testen\" bullet and no longer reflects the last\n#    edit location).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Locate \"ab 19:30\" and note exactly where it starts/ends.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"ab 19:30\"\n$found = $find.Execute()\n\nif ($found) {\n    $startPos = $rng.Start\n    $endPos = $rng.End\n\n    # Pin both edges of the match with bookmarks first so the surrounding\n    # \"(\" and \")\" runs are not swept up into the text replacement below; the\n    # trailing pin doubles as the new \"_GoBack\" location.\n    $d.Bookmarks.Add(\"ZZZTempPin\", $d.Range($startPos, $startPos))\n    $d.Bookmarks.Add(\"_GoBack\", $d.Range($endPos, $endPos))\n\n    # 3) Replace the time text itself.\n    $timeRange = $d.Range($startPos, $endPos)\n    $timeRange.Text = \"ab 19:15\"\n\n    # Remove the temporary helper bookmark again.\n    if ($d.Bookmarks.Exists(\"ZZZTempPin\")) {\n        $d.Bookmarks.Item(\"ZZZTempPin\").Delete()\n    }\n}\n"}
